# Apply the "update scripts wuth new tpm" edit:
#  - remove the 5 rows whose Sending cluster is "Resolving-Mac" (rows 22-26)
#  - refresh the recomputed TPM-derived statistics for the remaining rows
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Drop the "Resolving-Mac" sending-cluster block entirely (old rows 22:26).
$ws.Range("A22:A26").EntireRow.Delete()

# New TPM-derived values (columns E:T) for the remaining rows 2:21.
$rowData = @(
    @{Row=2; E=3; F=1; G=20.25895366666667; H=60.776861; I=0.9735054450004206; J=0.9760828245357948; K=3; L=1; M=1.400501333333333; N=4.201504; O=0.00926314904242919; P=0.009687730200823723; Q=28.37269162210489; R=255.354224598944; S=0.009017726030655249; T=0.009456027057760741},
    @{Row=3; E=3; F=1; G=20.25895366666667; H=60.776861; I=0.9735054450004206; J=0.9760828245357948; K=3; L=1; M=21.254561; N=63.763683; O=0.1405812059498714; P=0.1470248171880475; Q=430.5951665043403; R=3875.356498539063; S=0.1368565694569253; T=0.1435083988377682},
    @{Row=4; E=3; F=1; G=20.25895366666667; H=60.776861; I=0.9735054450004206; J=0.9760828245357948; K=3; L=1; M=63.87756733333333; N=191.632702; O=0.4224968677952986; P=0.4418622271050682; Q=1294.09267694538; R=11646.83409250842; S=0.411303001294346; T=0.4312941306883918},
    @{Row=5; E=3; F=1; G=20.25895366666667; H=60.776861; I=0.9735054450004206; J=0.9760828245357948; K=2; L=1; M=19.878555; N=39.75711; O=0.1314800731212866; P=0.0916710195312133; Q=402.718724705285; R=2416.31234823171; S=0.1279965670926259; T=0.08947850767210269},
    @{Row=6; E=3; F=1; G=20.25895366666667; H=60.776861; I=0.9735054450004206; J=0.9760828245357948; K=3; L=1; M=44.77944466666667; N=134.338334; O=0.2961787040911142; P=0.3097542059748472; Q=907.1846947210638; R=8164.662252489574; S=0.288331581125868; T=0.3023457602797712},
    @{Row=7; E=2; F=0.6666666666666666; G=0.1459843333333333; H=0.437953; I=0.007014999181255333; J=0.007033571563591034; K=3; L=1; M=1.400501333333333; N=4.201504; O=0.00926314904242919; P=0.009687730200823723; Q=0.2044512534791111; R=1.840061281312; S=0.00006498098294848689; T=0.00006813934365625579},
    @{Row=8; E=2; F=0.6666666666666666; G=0.1459843333333333; H=0.437953; I=0.007014999181255333; J=0.007033571563591034; K=3; L=1; M=21.254561; N=63.763683; O=0.1405812059498714; P=0.1470248171880475; Q=3.102832917877667; R=27.925496260899; S=0.0009861770446382351; T=0.001034109573316021},
    @{Row=9; E=2; F=0.6666666666666666; G=0.1459843333333333; H=0.437953; I=0.007014999181255333; J=0.007033571563591034; K=3; L=1; M=63.87756733333333; N=191.632702; O=0.4224968677952986; P=0.4418622271050682; Q=9.325124082111778; R=83.92611673900599; S=0.002963815181666962; T=0.003107869595591211},
    @{Row=10; E=2; F=0.6666666666666666; G=0.1459843333333333; H=0.437953; I=0.007014999181255333; J=0.007033571563591034; K=2; L=1; M=19.878555; N=39.75711; O=0.1314800731212866; P=0.0916710195312133; Q=2.901957599305; R=17.41174559583; S=0.0009223326052972165; T=0.0006447746761801402},
    @{Row=11; E=2; F=0.6666666666666666; G=0.1459843333333333; H=0.437953; I=0.007014999181255333; J=0.007033571563591034; K=3; L=1; M=44.77944466666667; N=134.338334; O=0.2961787040911142; P=0.3097542059748472; Q=6.537097376700222; R=58.833876390302; S=0.002077693366704432; T=0.002178678374847405},
    @{Row=12; E=2; F=0.6666666666666666; G=0.2405246666666667; H=0.721574; I=0.0115579548929112; J=0.0115885548618839; K=3; L=1; M=1.400501333333333; N=4.201504; O=0.00926314904242919; P=0.009687730200823723; Q=0.3368551163662223; R=3.031696047296; S=0.0001070630587987101; T=0.0001122667929193752},
    @{Row=13; E=2; F=0.6666666666666666; G=0.2405246666666667; H=0.721574; I=0.0115579548929112; J=0.0115885548618839; K=3; L=1; M=21.254561; N=63.763683; O=0.1405812059498714; P=0.1470248171880475; Q=5.112246199671334; R=46.010215797042; S=0.001624831237159673; T=0.001703805160042138},
    @{Row=14; E=2; F=0.6666666666666666; G=0.2405246666666667; H=0.721574; I=0.0115579548929112; J=0.0115885548618839; K=3; L=1; M=63.87756733333333; N=191.632702; O=0.4224968677952986; P=0.4418622271050682; Q=15.36413059032756; R=138.277175312948; S=0.004883199740374326; T=0.005120544660201285},
    @{Row=15; E=2; F=0.6666666666666666; G=0.2405246666666667; H=0.721574; I=0.0115579548929112; J=0.0115885548618839; K=2; L=1; M=19.878555; N=39.75711; O=0.1314800731212866; P=0.0916710195312133; Q=4.78128281519; R=28.68769689114; S=0.001519640754452496; T=0.001062334639082295},
    @{Row=16; E=2; F=0.6666666666666666; G=0.2405246666666667; H=0.721574; I=0.0115579548929112; J=0.0115885548618839; K=3; L=1; M=44.77944466666667; N=134.338334; O=0.2961787040911142; P=0.3097542059748472; Q=10.77056100196845; R=96.93504901771601; S=0.003423220102125991; T=0.003589603609638801},
    @{Row=17; E=1; F=0.5; G=0.164851; H=0.329702; I=0.007921600925412929; J=0.005295049038730392; K=3; L=1; M=1.400501333333333; N=4.201504; O=0.00926314904242919; P=0.009687730200823723; Q=0.2308740453013333; R=1.385244271808; S=0.00007337897002674496; T=0.00005129700648735103},
    @{Row=18; E=1; F=0.5; G=0.164851; H=0.329702; I=0.007921600925412929; J=0.005295049038730392; K=3; L=1; M=21.254561; N=63.763683; O=0.1405812059498714; P=0.1470248171880475; Q=3.503835635411; R=21.023013812466; S=0.001113628211148167; T=0.0007785036169210823},
    @{Row=19; E=1; F=0.5; G=0.164851; H=0.329702; I=0.007921600925412929; J=0.005295049038730392; K=3; L=1; M=63.87756733333333; N=191.632702; O=0.4224968677952986; P=0.4418622271050682; Q=10.53028085246733; R=63.181685114804; S=0.003346851578911301; T=0.002339682160883962},
    @{Row=20; E=1; F=0.5; G=0.164851; H=0.329702; I=0.007921600925412929; J=0.005295049038730392; K=2; L=1; M=19.878555; N=39.75711; O=0.1314800731212866; P=0.0916710195312133; Q=3.276999670305; R=13.10799868122; S=0.001041532668910943; T=0.000485402543848186},
    @{Row=21; E=1; F=0.5; G=0.164851; H=0.329702; I=0.007921600925412929; J=0.005295049038730392; K=3; L=1; M=44.77944466666667; N=134.338334; O=0.2961787040911142; P=0.3097542059748472; Q=7.381936232744667; R=44.291617396468; S=0.002346209496415773; T=0.00164016371058981}
)

foreach ($r in $rowData) {
    $ws.Cells.Item($r.Row, 5).Value = $r.E
    $ws.Cells.Item($r.Row, 6).Value = $r.F
    $ws.Cells.Item($r.Row, 7).Value = $r.G
    $ws.Cells.Item($r.Row, 8).Value = $r.H
    $ws.Cells.Item($r.Row, 9).Value = $r.I
    $ws.Cells.Item($r.Row, 10).Value = $r.J
    $ws.Cells.Item($r.Row, 11).Value = $r.K
    $ws.Cells.Item($r.Row, 12).Value = $r.L
    $ws.Cells.Item($r.Row, 13).Value = $r.M
    $ws.Cells.Item($r.Row, 14).Value = $r.N
    $ws.Cells.Item($r.Row, 15).Value = $r.O
    $ws.Cells.Item($r.Row, 16).Value = $r.P
    $ws.Cells.Item($r.Row, 17).Value = $r.Q
    $ws.Cells.Item($r.Row, 18).Value = $r.R
    $ws.Cells.Item($r.Row, 19).Value = $r.S
    $ws.Cells.Item($r.Row, 20).Value = $r.T
}
